$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking price strings
# (e.g. "0.9986") are stored as literal text, matching the source data,
# instead of being auto-coerced into numbers by the input parser.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '31.301.89'
$ws.Range('E2').Value = '  +1.85%  '
$ws.Range('D3').Value = '2.002.02'
$ws.Range('E3').Value = '  +5.44%  '
$ws.Range('D4').Value = '0.9986'
$ws.Range('E4').Value = '  +0.39%  '
$ws.Range('D5').Value = '0.7789'
$ws.Range('E5').Value = '  +43.36%  '
$ws.Range('D6').Value = '255.56'
$ws.Range('E6').Value = '  +3.46%  '
$ws.Range('D7').Value = '0.9986'
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').Value = '0.3483'
$ws.Range('E8').Value = '  +18.75%  '
$ws.Range('D9').Value = '28.25'
$ws.Range('E9').Value = '  +24.45%  '
$ws.Range('D10').Value = '0.07187'
$ws.Range('E10').Value = '  +9.37%  '
$ws.Range('D11').Value = '0.8466'
$ws.Range('E11').Value = '  +10.49%  '
$ws.Range('D12').Value = '0.08212'
$ws.Range('E12').Value = '  +4.63%  '
$ws.Range('D13').Value = '101.26'
$ws.Range('E13').Value = '  +0.78%  '
$ws.Range('D14').Value = '5.671'
$ws.Range('E14').Value = '  +7.43%  '
$ws.Range('D15').Value = '2.001.53'
$ws.Range('E15').Value = '  +5.58%  '
$ws.Range('D16').Value = '15.44'
$ws.Range('E16').Value = '  +16.35%  '
$ws.Range('D17').Value = '273.31'
$ws.Range('E17').Value = '  -4.22%  '
$ws.Range('D18').Value = '31.288.73'
$ws.Range('E18').Value = '  +2.10%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.000008348'
$ws.Range('E19').Value = '  +10.28%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = '6.022'
$ws.Range('E20').Value = '  +11.81%  '
$ws.Range('D21').Value = '2.261.08'
$ws.Range('E21').Value = '  +6.35%  '
$ws.Range('D22').Value = '0.9982'
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('D23').Value = '0.9985'
$ws.Range('E23').Value = '  +0.33%  '
$ws.Range('D24').Value = '7.146'
$ws.Range('E24').Value = '  +10.39%  '
$ws.Range('D25').Value = '10.14'
$ws.Range('E25').Value = '  +9.79%  '
$ws.Range('D26').Value = '164.50'
$ws.Range('E26').Value = '  +0.37%  '
$ws.Range('D27').Value = '0.1425'
$ws.Range('E27').Value = '  +40.84%  '
$ws.Range('D28').Value = '20.04'
$ws.Range('E28').Value = '  +4.21%  '
$ws.Range('D29').Value = '2.405'
$ws.Range('E29').Value = '  +24.76%  '
$ws.Range('D30').Value = '1.599'
$ws.Range('E30').Value = '  +6.07%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '4.668'
$ws.Range('E31').Value = '  +8.70%  '
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').Value = '1.368'
$ws.Range('E32').Value = '  +1.92%  '
$ws.Range('D33').Value = '4.475'
$ws.Range('E33').Value = '  +6.08%  '
$ws.Range('D34').Value = '0.05368'
$ws.Range('E34').Value = '  +9.82%  '
$ws.Range('D35').Value = '1.277'
$ws.Range('E35').Value = '  +11.88%  '
$ws.Range('D36').Value = '0.7918'
$ws.Range('E36').Value = '  +12.23%  '
$ws.Range('D37').Value = '2.770'
$ws.Range('E37').Value = '  -0.40%  '
$ws.Range('D38').Value = '0.9980'
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('D39').Value = '0.02014'
$ws.Range('E39').Value = '  +5.01%  '
$ws.Range('D40').Value = '87.02'
$ws.Range('E40').Value = '  +14.35%  '
$ws.Range('D41').Value = '2.930'
$ws.Range('E41').Value = '  +1.85%  '
$ws.Range('D42').Value = '6.783'
$ws.Range('E42').Value = '  +6.99%  '
$ws.Range('D43').Value = '2.159'
$ws.Range('E43').Value = '  +8.46%  '
$ws.Range('D44').Value = '0.4714'
$ws.Range('E44').Value = '  +9.64%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').Value = '0.8593'
$ws.Range('E45').Value = '  +1.76%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').Value = '105.58'
$ws.Range('E46').Value = '  +4.04%  '
$ws.Range('D47').Value = '10.30'
$ws.Range('E47').Value = '  +3.09%  '
$ws.Range('B48').Value = 'SynthetixNetwork'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D48').Value = '3.196'
$ws.Range('E48').Value = '  +53.04%  '
$ws.Range('D49').Value = '0.9988'
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('D50').Value = '7.796'
$ws.Range('E50').Value = '  +9.75%  '
$ws.Range('D51').Value = '37.94'
$ws.Range('E51').Value = '  +7.18%  '

# Restore the default (General) display format / style so the saved
# cells look identical to the originals (no stray explicit style index).
$ws.Range("D2:D51").NumberFormat = "General"
$ws.Range("D2:D51").Style = "Normal"
